$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text interpretation for D and E columns so numeric-looking strings
# (e.g. "516.16") are stored as text, matching the original inlineStr cell type.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "68.755.71"
$ws.Range("E2").Value = "  -0.69%  "

$ws.Range("D3").Value = "3.849.93"
$ws.Range("E3").Value = "  -2.25%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "516.16"
$ws.Range("E5").Value = "  +4.51%  "

$ws.Range("D6").Value = "140.89"
$ws.Range("E6").Value = "  -4.22%  "

$ws.Range("D7").Value = "0.607"
$ws.Range("E7").Value = "  -2.70%  "

$ws.Range("E8").Value = "  +0.10%  "

$ws.Range("D9").Value = "0.713"
$ws.Range("E9").Value = "  -2.55%  "

$ws.Range("D10").Value = "0.168"
$ws.Range("E10").Value = "  -4.98%  "

$ws.Range("D11").Value = "0.0000322"
$ws.Range("E11").Value = "  -8.58%  "

$ws.Range("D12").Value = "41.50"
$ws.Range("E12").Value = "  -4.20%  "

$ws.Range("D13").Value = "10.28"
$ws.Range("E13").Value = "  -1.51%  "

$ws.Range("D14").Value = "4.456.75"
$ws.Range("E14").Value = "  -2.45%  "

$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "3.855.67"
$ws.Range("E15").Value = "  -1.92%  "

$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "21.09"
$ws.Range("E16").Value = "  +5.99%  "

$ws.Range("D17").Value = "13.95"
$ws.Range("E17").Value = "  -2.46%  "

$ws.Range("E18").Value = "  -2.20%  "

$ws.Range("E19").Value = "  +0.89%  "

$ws.Range("D20").Value = "68.596.07"
$ws.Range("E20").Value = "  -1.02%  "

$ws.Range("D21").Value = "414.76"
$ws.Range("E21").Value = "  -5.50%  "

$ws.Range("E22").Value = "  -0.54%  "

$ws.Range("D23").Value = "12.06"
$ws.Range("E23").Value = "  -0.39%  "

$ws.Range("D24").Value = "13.90"
$ws.Range("E24").Value = "  -4.37%  "

$ws.Range("D25").Value = "86.68"
$ws.Range("E25").Value = "  -2.62%  "

$ws.Range("E26").Value = "  +4.71%  "

$ws.Range("D27").Value = "10.38"
$ws.Range("E27").Value = "  -6.71%  "

$ws.Range("D28").Value = "35.37"
$ws.Range("E28").Value = "  -4.79%  "

$ws.Range("D29").Value = "13.34"
$ws.Range("E29").Value = "  -0.25%  "

$ws.Range("D30").Value = "676.61"
$ws.Range("E30").Value = "  -3.73%  "

$ws.Range("D31").Value = "6.98"
$ws.Range("E31").Value = "  +15.11%  "

$ws.Range("E32").Value = "  -1.74%  "

$ws.Range("D33").Value = "0.124"
$ws.Range("E33").Value = "  -5.09%  "

$ws.Range("D34").Value = "66.13"
$ws.Range("E34").Value = "  +6.91%  "

$ws.Range("D35").Value = "0.444"
$ws.Range("E35").Value = "  -4.63%  "

$ws.Range("D36").Value = "0.0₃0846"
$ws.Range("E36").Value = "  -5.81%  "

$ws.Range("D37").Value = "39.10"
$ws.Range("E37").Value = "  -4.17%  "

$ws.Range("D38").Value = "3.38"
$ws.Range("E38").Value = "  +10.10%  "

$ws.Range("E39").Value = "  -3.31%  "

$ws.Range("E40").Value = "  +0.29%  "

$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  -0.24%  "

$ws.Range("D42").Value = "0.0473"
$ws.Range("E42").Value = "  -3.12%  "

$ws.Range("E43").Value = "  +4.66%  "

$ws.Range("D44").Value = "2.84"
$ws.Range("E44").Value = "  -1.98%  "

$ws.Range("D45").Value = "3.60"
$ws.Range("E45").Value = "  +6.22%  "

$ws.Range("E46").Value = "  -1.44%  "

$ws.Range("D47").Value = "0.000286"
$ws.Range("E47").Value = "  +19.02%  "

$ws.Range("D48").Value = "3.02"
$ws.Range("E48").Value = "  +0.18%  "

$ws.Range("D49").Value = "3.29"
$ws.Range("E49").Value = "  -2.57%  "

$ws.Range("D50").Value = "143.05"
$ws.Range("E50").Value = "  -0.88%  "

$ws.Range("D51").Value = "8.69"
$ws.Range("E51").Value = "  +1.36%  "

# Restore the default (Normal) style on the data range so no stray
# quote-prefix / text-format styling remains applied to the cells.
$dataRange.Style = "Normal"
